# "download articles with pandoc title blocks"
#
# The heading paragraph ("On Pilgrimage - September 1946", styled
# Heading1 and wrapped in a bookmark) and the byline paragraph
# ("By Dorothy Day", bold) are replaced by a pandoc-style title block:
#   - a Title-styled paragraph with the article title, split word-by-word
#     (and space-by-space) into separate runs
#   - an Authors-styled paragraph with just "Dorothy Day" (the leading
#     "By " is dropped), likewise split into separate runs
#
# Both the "Title" and "Authors" paragraph styles already exist in
# styles.xml, so only document.xml content changes.

$d = $word.ActiveDocument

$p1 = $d.Paragraphs.Item(1)
$p2 = $d.Paragraphs.Item(2)

# Range spanning both the heading paragraph and the byline paragraph
# (including both of their end-of-paragraph marks), so the new
# paragraph styles actually take effect on save.
$target = $d.Range($p1.Range.Start, $p2.Range.End)

$newPkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
'<w:p><w:pPr><w:pStyle w:val="Title"/></w:pPr>' + `
'<w:r><w:t xml:space="preserve">On</w:t></w:r>' + `
'<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
'<w:r><w:t xml:space="preserve">Pilgrimage</w:t></w:r>' + `
'<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
'<w:r><w:t xml:space="preserve">-</w:t></w:r>' + `
'<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
'<w:r><w:t xml:space="preserve">September</w:t></w:r>' + `
'<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
'<w:r><w:t xml:space="preserve">1946</w:t></w:r>' + `
'</w:p>' + `
'<w:p><w:pPr><w:pStyle w:val="Authors"/></w:pPr>' + `
'<w:r><w:t xml:space="preserve">Dorothy</w:t></w:r>' + `
'<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
'<w:r><w:t xml:space="preserve">Day</w:t></w:r>' + `
'</w:p>' + `
'</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($newPkg) | Out-Null

Write-Output ("ParaCount=" + $d.Paragraphs.Count)
Write-Output ("p1 style=" + $d.Paragraphs.Item(1).Style.NameLocal + " text=[" + $d.Paragraphs.Item(1).Range.Text + "]")
Write-Output ("p2 style=" + $d.Paragraphs.Item(2).Style.NameLocal + " text=[" + $d.Paragraphs.Item(2).Range.Text + "]")
